$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("RF")
$ws.Range("C2").Value = 0.6807784467545627
$ws.Range("D2").Value = 0.02553605817873348
$ws.Range("C3").Value = 0.8321025488740412
$ws.Range("D3").Value = 0.01749523915829225
$ws.Range("C4").Value = 0.8309455519383632
$ws.Range("D4").Value = 0.02481878422344192
$ws.Range("C5").Value = 0.8548101523078547
$ws.Range("D5").Value = 0.02397779369039548
$ws.Range("C6").Value = 0.807194
$ws.Range("D6").Value = 0.02830178771305979
$ws.Range("C7").Value = 0.8423684334437719
$ws.Range("D7").Value = 0.01751306892569684
$ws.Range("C8").Value = 0.831918886446506
$ws.Range("D8").Value = 0.01752696332899967
$ws.Range("C9").Value = 0.8311479930248875
$ws.Range("D9").Value = 0.01754911812571582
$ws.Range("C10").Value = 0.8310053743143534
$ws.Range("D10").Value = 0.01745472697808747
$ws.Range("C11").Value = 0.6634534029714765
$ws.Range("D11").Value = 0.03491329361663161
$ws.Range("C12").Value = 0.833958
$ws.Range("D12").Value = 0.02615241893834164
$ws.Range("C13").Value = 0.8310053743143534
$ws.Range("D13").Value = 0.01745472697808745

$ws = $wb.Worksheets.Item("LGBM")
$ws.Range("C2").Value = 0.7025111363086483
$ws.Range("D2").Value = 0.02470783379758652
$ws.Range("C3").Value = 0.829962583518931
$ws.Range("D3").Value = 0.01679499426843057
$ws.Range("C4").Value = 0.841218390448103
$ws.Range("D4").Value = 0.02327805787531907
$ws.Range("C5").Value = 0.8341142869635675
$ws.Range("D5").Value = 0.02649700475678737
$ws.Range("C6").Value = 0.8254319999999998
$ws.Range("D6").Value = 0.02640871408440473
$ws.Range("C7").Value = 0.8372960572348511
$ws.Range("D7").Value = 0.01785623434010699
$ws.Range("C8").Value = 0.8299682626625777
$ws.Range("D8").Value = 0.01678779859508431
$ws.Range("C9").Value = 0.8293487036261717
$ws.Range("D9").Value = 0.01675474004540455
$ws.Range("C10").Value = 0.8297723184136572
$ws.Range("D10").Value = 0.01672471598967783
$ws.Range("C11").Value = 0.6595099890663465
$ws.Range("D11").Value = 0.03348126304242634
$ws.Range("C12").Value = 0.8182579999999998
$ws.Range("D12").Value = 0.02592201089768291
$ws.Range("C13").Value = 0.8297723184136572
$ws.Range("D13").Value = 0.01672471598967783

$ws = $wb.Worksheets.Item("XGB")
$ws.Range("C2").Value = 0.7125768740230027
$ws.Range("D2").Value = 0.02391794857403803
$ws.Range("C3").Value = 0.8366891363523883
$ws.Range("D3").Value = 0.0157735359493348
$ws.Range("C4").Value = 0.8406516858264177
$ws.Range("D4").Value = 0.02242608470366684
$ws.Range("C5").Value = 0.8508306644282639
$ws.Range("D5").Value = 0.02467303713115537
$ws.Range("C6").Value = 0.8210400000000001
$ws.Range("D6").Value = 0.02646812484202577
$ws.Range("C7").Value = 0.8453755927391707
$ws.Range("D7").Value = 0.01625791113675387
$ws.Range("C8").Value = 0.8366000091863516
$ws.Range("D8").Value = 0.01576197674870596
$ws.Range("C9").Value = 0.835903497881561
$ws.Range("D9").Value = 0.01578800859274338
$ws.Range("C10").Value = 0.8359347923833752
$ws.Range("D10").Value = 0.01574509591048038
$ws.Range("C11").Value = 0.6726286358788031
$ws.Range("D11").Value = 0.03166476891507117
$ws.Range("C12").Value = 0.832742
$ws.Range("D12").Value = 0.02526380883268976
$ws.Range("C13").Value = 0.8359347923833752
$ws.Range("D13").Value = 0.01574509591048038

$ws = $wb.Worksheets.Item("KNN")
$ws.Range("C2").Value = 0.6603189337684975
$ws.Range("D2").Value = 0.03320883754884141
$ws.Range("C3").Value = 0.8202138084632519
$ws.Range("D3").Value = 0.01899122848294163
$ws.Range("C4").Value = 0.8296144686962088
$ws.Range("D4").Value = 0.02107775356168779
$ws.Range("C5").Value = 0.8282887685886745
$ws.Range("D5").Value = 0.03143914656983188
$ws.Range("C6").Value = 0.81145
$ws.Range("D6").Value = 0.02428564915957655
$ws.Range("C7").Value = 0.8285593772450263
$ws.Range("D7").Value = 0.01949311707121323
$ws.Range("C8").Value = 0.8202043930832749
$ws.Range("D8").Value = 0.01891162044003608
$ws.Range("C9").Value = 0.8195222143366143
$ws.Range("D9").Value = 0.01900188818496525
$ws.Range("C10").Value = 0.8198705913914577
$ws.Range("D10").Value = 0.01880093374156847
$ws.Range("C11").Value = 0.6399082627030944
$ws.Range("D11").Value = 0.03822011155659995
$ws.Range("C12").Value = 0.8104600000000001
$ws.Range("D12").Value = 0.03271734234004832
$ws.Range("C13").Value = 0.8198705913914577
$ws.Range("D13").Value = 0.01880093374156849

$ws = $wb.Worksheets.Item("SVM")
$ws.Range("C2").Value = 0.7134828397952719
$ws.Range("D2").Value = 0.02501097003083421
$ws.Range("C3").Value = 0.8399814897302649
$ws.Range("D3").Value = 0.01763945936856931
$ws.Range("C4").Value = 0.8468281354585458
$ws.Range("D4").Value = 0.02340900591933882
$ws.Range("C5").Value = 0.8498359181225917
$ws.Range("D5").Value = 0.02673675081207115
$ws.Range("C6").Value = 0.829298
$ws.Range("D6").Value = 0.02840155564868422
$ws.Range("C7").Value = 0.8479380539505488
$ws.Range("D7").Value = 0.01702802179071703
$ws.Range("C8").Value = 0.8399422566386806
$ws.Range("D8").Value = 0.01764995142822005
$ws.Range("C9").Value = 0.8393027725546232
$ws.Range("D9").Value = 0.01779203305965085
$ws.Range("C10").Value = 0.839566994870975
$ws.Range("D10").Value = 0.01782901681909002
$ws.Range("C11").Value = 0.6794708419547071
$ws.Range("D11").Value = 0.03571194622536493
$ws.Range("C12").Value = 0.8329840000000001
$ws.Range("D12").Value = 0.02892617175232315
$ws.Range("C13").Value = 0.839566994870975
$ws.Range("D13").Value = 0.01782901681909003
